# Reorders the "Recorded By" (column G) comma-separated list of names so that
# any entry containing "system" (case-insensitive) is moved to the front of
# the list (preserving relative order), followed by the remaining entries
# (e.g. email addresses), also preserving their relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2

    if ($current -ne $null -and $current -ne "") {
        $parts = $current -split ","
        $systemParts = @()
        $otherParts = @()

        foreach ($part in $parts) {
            $trimmed = $part.Trim()
            if ($trimmed -match "(?i)system") {
                $systemParts += $trimmed
            } else {
                $otherParts += $trimmed
            }
        }

        $newValue = ($systemParts + $otherParts) -join ", "

        if ($newValue -ne $current) {
            $cell.Value2 = $newValue
        }
    }
}
